$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.22
$ws.Range("G2").Value = 1.24
$ws.Range("H2").Value = 25
$ws.Range("I2").Value = 950
$ws.Range("J2").Value = 5.9
$ws.Range("K2").Value = 6.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2.2
$ws.Range("O2").Value = 1.75
$ws.Range("P2").Value = 1.2
$ws.Range("Q2").Value = 5.7
$ws.Range("R2").Value = 1.02
$ws.Range("S2").Value = 26
$ws.Range("T2").Value = 3.05
$ws.Range("U2").Value = 1.26
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 5.2
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 2.46
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 690
$ws.Range("AG2").Value = 990
$ws.Range("AH2").Value = 990
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 990
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000
